# abnormal.xlsx - "enemy" target removed from abnormal-data rows.
# The "target" column (E) for the Slow/Snare abnormal entries used to read
# "Enemy"; per the commit ("enemy 삭제") it is changed to "target" for those
# two rows. The "faction" targeted row (E4) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "target"
$ws.Range("E3").Value = "target"

# Reflect the author's final on-screen selection (cosmetic, no data impact).
$ws.Range("J13:K13").Select()
